$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 13.33
$ws.Range("E4").Value = 12.919
$ws.Range("E7").Value = 13.35
$ws.Range("E8").Value = 12.913
$ws.Range("B11").Value = 6.313
$ws.Range("B12").Value = 5.915
$ws.Range("E12").Value = 13.143
$ws.Range("E14").Value = 12.911
$ws.Range("B15").Value = 6.336
$ws.Range("E22").Value = 12.862
